# Auto-generated edit script applying cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.493.67"
$ws.Range("E2").Value = "  +3.83%  "
$ws.Range("D3").Value = "2.622.18"
$ws.Range("E3").Value = "  +1.85%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.93%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("E8").Value = "  +4.24%  "
$ws.Range("D9").Value = "2.639.60"
$ws.Range("E9").Value = "  +2.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.74"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.76%  "
$ws.Range("E11").Value = "  +5.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.153"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +11.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.345"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.18%  "
$ws.Range("D14").Value = "3.080.37"
$ws.Range("E14").Value = "  +1.35%  "
$ws.Range("D15").Value = "60.444.75"
$ws.Range("E15").Value = "  +3.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000139"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.75%  "
$ws.Range("D18").Value = "2.631.76"
$ws.Range("E18").Value = "  +2.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "342.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.51%  "
$ws.Range("E21").Value = "  +4.04%  "
$ws.Range("E22").Value = "  +4.19%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("E25").Value = "  +7.80%  "
$ws.Range("E26").Value = "  +3.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.01%  "
$ws.Range("D29").Value = "0.0₃0800"
$ws.Range("E29").Value = "  +9.28%  "
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("E31").Value = "  +4.96%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "160.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.13%  "
$ws.Range("E34").Value = "  +2.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.12"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.34%  "
$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.901"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.57%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.16"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.887"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.32%  "
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.52"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "298.67"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.76%  "
$ws.Range("E43").Value = "  -0.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0983"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.602"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0543"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.48%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0237"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.28%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "126.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +16.09%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.22%  "
